# completed extent report with screenshot
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Add a new row 4 reusing the same values/styles as row 2/3 (Username / Admin12345@gmail.com, password / Admin@123)
$ws.Range("A2:B2").Copy($ws.Range("A4:B4"))

# Move the active selection to A5, below the newly added row
$ws.Range("A5").Select()
